# Update the "想去人数" (want-to-go count) figures that were refreshed in the
# gh-pages data regeneration (commit 456a3b4).
#
# Sheet "展览" (Exhibitions) - column F
#   F4  962  -> 963
#   F6  7406 -> 7410
#   F9  6639 -> 6641
#   F11 283  -> 284
#   F12 4636 -> 4637
#   F16 4791 -> 4794
#   F17 4791 -> 4794
#   F31 8452 -> 8454
#   F47 861  -> 862
#   F48 1155 -> 1156
#
# Sheet "全部类型" (All types) - column F (same events, different row numbers
# because this sheet aggregates rows from every category)
#   F6  962  -> 963
#   F8  7406 -> 7410
#   F11 6639 -> 6641
#   F13 283  -> 284
#   F15 4636 -> 4637
#   F19 4791 -> 4794
#   F33 8452 -> 8454
#   F48 861  -> 862
#   F49 1155 -> 1156

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    "F4"  = 963
    "F6"  = 7410
    "F9"  = 6641
    "F11" = 284
    "F12" = 4637
    "F16" = 4794
    "F17" = 4794
    "F31" = 8454
    "F47" = 862
    "F48" = 1156
}
foreach ($cell in $exhibitionUpdates.Keys) {
    $wsExhibitions.Range($cell).Value = $exhibitionUpdates[$cell]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    "F6"  = 963
    "F8"  = 7410
    "F11" = 6641
    "F13" = 284
    "F15" = 4637
    "F19" = 4794
    "F33" = 8454
    "F48" = 862
    "F49" = 1156
}
foreach ($cell in $allTypesUpdates.Keys) {
    $wsAllTypes.Range($cell).Value = $allTypesUpdates[$cell]
}
